# Fruta / hortaliza, semanal
# Rotates the daily price records among rows 2-5 and 7-10 (row 6 is unchanged)
# as described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each row, taken from the previous state of another row
# (columns D, J, K, L, M, N, O, P, Q). All other columns are identical
# across rows and remain untouched.

$rows = @{
    2  = @{ D = 44215; J = 140; K = 3500; L = 4000; M = 3768; N = "`$/paquete 2 kilos";   O = "Provincia de Diguillín"; P = 1884; Q = 2 }
    3  = @{ D = 44161; J = 50;  K = 2800; L = 3000; M = 2900; N = "`$/paquete 2 kilos";   O = "Provincia de Diguillín"; P = 1450; Q = 2 }
    4  = @{ D = 44209; J = 150; K = 3500; L = 4000; M = 3767; N = "`$/paquete 2 kilos";   O = "Provincia de Diguillín"; P = 1884; Q = 2 }
    5  = @{ D = 44210; J = 105; K = 3500; L = 4000; M = 3714; N = "`$/paquete 2 kilos";   O = "Provincia de Diguillín"; P = 1857; Q = 2 }
    7  = @{ D = 44223; J = 80;  K = 3500; L = 3800; M = 3688; N = "`$/paquete 2 kilos";   O = "Provincia de Diguillín"; P = 1844; Q = 2 }
    8  = @{ D = 44225; J = 80;  K = 3400; L = 3700; M = 3550; N = "`$/paquete 2 kilos";   O = "Provincia de Diguillín"; P = 1775; Q = 2 }
    9  = @{ D = 44208; J = 85;  K = 3700; L = 4000; M = 3824; N = "`$/paquete 2 kilos";   O = "Provincia de Diguillín"; P = 1912; Q = 2 }
    10 = @{ D = 44160; J = 43;  K = 3500; L = 4000; M = 3709; N = "`$/paquete 36 unidades"; O = "Región Metropolitana";   P = 103;  Q = 36 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $vals.N   # N: Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $vals.O   # O: Origen
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q: Kg o Unidades
}
